$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '41.763.58'
Set-TextValue 'E2' '  -1.32%  '
Set-TextValue 'D3' '2.221.77'
Set-TextValue 'E3' '  -0.98%  '
Set-TextValue 'D5' '251.10'
Set-TextValue 'E5' '  +6.01%  '
Set-TextValue 'D6' '0.631'
Set-TextValue 'E6' '  -0.14%  '
Set-TextValue 'D7' '71.65'
Set-TextValue 'E7' '  +3.06%  '
Set-TextValue 'E8' '  -0.14%  '
Set-TextValue 'D9' '0.604'
Set-TextValue 'E9' '  +8.78%  '
Set-TextValue 'D10' '40.50'
Set-TextValue 'E10' '  +10.12%  '
Set-TextValue 'D11' '0.0964'
Set-TextValue 'E11' '  -2.84%  '
Set-TextValue 'D12' '58.38'
Set-TextValue 'E12' '  -0.82%  '
Set-TextValue 'D13' '7.23'
Set-TextValue 'E13' '  +7.34%  '
Set-TextValue 'D14' '0.106'
Set-TextValue 'E14' '  -0.58%  '
Set-TextValue 'D15' '2.552.70'
Set-TextValue 'E15' '  -1.20%  '
Set-TextValue 'D16' '14.97'
Set-TextValue 'E16' '  -0.49%  '
Set-TextValue 'E17' '  +1.11%  '
Set-TextValue 'D18' '2.212.96'
Set-TextValue 'E18' '  -1.69%  '
Set-TextValue 'D19' '41.713.91'
Set-TextValue 'E19' '  -1.31%  '
Set-TextValue 'D20' '0.0₃0963'
Set-TextValue 'E20' '  -1.32%  '
Set-TextValue 'D21' '6.22'
Set-TextValue 'E21' '  -0.61%  '
Set-TextValue 'D22' '72.87'
Set-TextValue 'E22' '  -0.69%  '
Set-TextValue 'D23' '232.49'
Set-TextValue 'E23' '  -1.64%  '
Set-TextValue 'E24' '  +5.22%  '
Set-TextValue 'E25' '  +9.33%  '
Set-TextValue 'D26' '1.00'
Set-TextValue 'E26' '  -0.04%  '
Set-TextValue 'D28' '10.87'
Set-TextValue 'E28' '  +8.59%  '
Set-TextValue 'D29' '170.82'
Set-TextValue 'E29' '  -0.16%  '
Set-TextValue 'E30' '  -6.42%  '
Set-TextValue 'D31' '20.79'
Set-TextValue 'E31' '  +1.09%  '
Set-TextValue 'D32' '0.123'
Set-TextValue 'B33' 'Stellar'
Set-TextValue 'C33' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D33' '0.125'
Set-TextValue 'E33' '  -1.34%  '
Set-TextValue 'B34' 'InternetComputer(DFINITY)'
Set-TextValue 'C34' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D34' '5.58'
Set-TextValue 'E34' '  +4.82%  '
Set-TextValue 'D35' '0.0740'
Set-TextValue 'E35' '  +2.83%  '
Set-TextValue 'E36' '  +0.64%  '
Set-TextValue 'D37' '26.37'
Set-TextValue 'E37' '  +17.72%  '
Set-TextValue 'D38' '4.03'
Set-TextValue 'E38' '  +8.78%  '
Set-TextValue 'E39' '  +10.07%  '
Set-TextValue 'E40' '  +0.26%  '
Set-TextValue 'D41' '5.95'
Set-TextValue 'E41' '  +0.39%  '
Set-TextValue 'B42' 'Celestia'
Set-TextValue 'C42' 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue 'D42' '12.31'
Set-TextValue 'E42' '  +21.25%  '
Set-TextValue 'B43' 'MultiversX'
Set-TextValue 'C43' 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue 'D43' '65.68'
Set-TextValue 'E43' '  +1.00%  '
Set-TextValue 'D44' '0.204'
Set-TextValue 'E44' '  +6.52%  '
Set-TextValue 'D45' '4.88'
Set-TextValue 'E45' '  -1.03%  '
Set-TextValue 'D46' '8.71'
Set-TextValue 'E46' '  -6.85%  '
Set-TextValue 'B47' 'SynthetixNetwork'
Set-TextValue 'C47' 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
Set-TextValue 'D47' '4.75'
Set-TextValue 'E47' '  +3.09%  '
Set-TextValue 'B48' 'Cronos'
Set-TextValue 'C48' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D48' '0.102'
Set-TextValue 'E48' '  -0.95%  '
Set-TextValue 'E49' '  -0.17%  '
Set-TextValue 'E50' '  +5.63%  '
Set-TextValue 'D51' '2.39'
Set-TextValue 'E51' '  +2.37%  '
